# Adding profits ("M_PL") column to the region table.
# A new column is inserted right after the existing "M_ETR" column (column B),
# pushing the remaining metric columns one slot to the right, and the new
# column is populated with a header label and per-row profit figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C - shifts old C:J to D:K
$ws.Range("C1").EntireColumn.Insert()

# Header for the newly inserted column
$ws.Range("C1").Value = "M_PL"

# Profit values for rows 2-8 (new column C)
$ws.Range("C2").Value = 106960829103
$ws.Range("C3").Value = 145933306887
$ws.Range("C4").Value = 45685889210
$ws.Range("C5").Value = 10524807277
$ws.Range("C6").Value = 885447038872
$ws.Range("C7").Value = 12956669707
$ws.Range("C8").Value = 5046999058
